# dev of account system
# Rework the GameCenter account-binding error strings: rename/re-split the
# existing keys+messages and append a new "not yet bound" error pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174: playerAlreadyBindGCAId -> userAlreadyBindGCAId (message unchanged)
$ws.Range("A174").Value = "userAlreadyBindGCAId"

# Row 175: theGCIdAlreadyHasDatas -> theGCIdAlreadyBindedByOtherUser
#          message: 此GameCenter账号已有玩家数据 -> 此GameCenter账号已被其他玩家绑定
$ws.Range("A175").Value = "theGCIdAlreadyBindedByOtherUser"
$ws.Range("C175").Value = "此GameCenter账号已被其他玩家绑定"

# Row 176: theGCAccountDoNotHasData -> theGCAccountIsNotBindedByOtherUser
#          message: 此GameCenter账号下无玩家数据 -> 此GameCenter账号未被其他玩家绑定
$ws.Range("A176").Value = "theGCAccountIsNotBindedByOtherUser"
$ws.Range("C176").Value = "此GameCenter账号未被其他玩家绑定"

# New row 177: theUserDoNotBindGCId / 676 / 当前玩家还未绑定GameCenter账号
# Seed the new row by copying row 176's formatting (and values), then
# overwrite with the real content so the style/height match the rest of
# the table.
$ws.Range("A176:C176").Copy($ws.Range("A177:C177")) | Out-Null
$ws.Rows.Item(177).RowHeight = $ws.Rows.Item(176).RowHeight

$ws.Range("A177").Value = "theUserDoNotBindGCId"
$ws.Range("B177").Value = 676
$ws.Range("C177").Value = "当前玩家还未绑定GameCenter账号"

# Leave the cursor one row below the newly-added row, matching where Excel
# lands after typing the last entry and pressing Enter.
$ws.Range("A178").Select() | Out-Null
